# Generate Report for Handoff
# Refresh the localization-status report: swap the superseded source file's
# GUID/hash for the new one everywhere it appears, stamp the new handoff
# timestamps, and clear out the (now stale) "Latest Target File" /
# "Latest Handback File" columns + their hyperlinks on the language sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "977574a9-6d7b-41ee-be0e-dfcb9cd2cb0a"
$newGuid = "d948397b-0740-46d8-89be-8a2807885567"
$oldHash = "bd29fc494909e1a9accfe9744ca33797666fa5cd"
$newHash = "38dea0824b1cd926ecc7ebbe74591623a731e12a"

# 1) Swap the GUID and content-hash everywhere (file names, paths, hyperlink
#    display text is handled separately below) across every sheet.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldGuid, $newGuid)
    $ws.Cells.Replace($oldHash, $newHash)
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# 2) Update the hyperlink display text (TextToDisplay) so it matches the new
#    file name too - Replace() only touches cell values, not hyperlink text.
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\" + $newGuid + ".md"
}
foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newGuid + ".md"
    }
}
foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newGuid + ".md"
    }
}

# 3) New handoff timestamps.
$wsOverview.Range("G2").Value = "2016-08-29 13:00:23"
$wsZhCn.Range("H2").Value = "2016-08-29 12:59:59"
$wsDeDe.Range("H2").Value = "2016-08-29 13:00:23"

# 4) The old handback info is no longer valid - clear "Latest Target File"
#    (I2) and "Latest Handback File" (J2) and reset the handback datetime
#    (K2) on both language sheets, dropping the now-dead I2 hyperlinks.
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$I$2') {
            $hl.Delete()
        }
    }
    $ws.Range("I2").Value = ""
    $ws.Range("I2").Style = "Normal"
    $ws.Range("J2").Value = ""
    $ws.Range("K2").Value = "0001-01-01 00:00:00"
}

# 5) The Latest Target File / Latest Handback File columns are now short
#    blanks instead of full file names - narrow them back down to fit.
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(9).ColumnWidth = 17.8
    $ws.Columns.Item(10).ColumnWidth = 20.8
}
